# Update the "dSF" column (F) values for specific rows per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    3  = 5
    12 = 2
    16 = 1
    21 = 3
    24 = 3
    25 = 3
    38 = 1
    40 = 3
    46 = 0
    48 = 0
    51 = 3
    58 = 0
    62 = 3
    66 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
